$d = $word.ActiveDocument

# Smart quotes used verbatim in the "Git commit" bullet's text.
$ldq = [char]0x201C
$rdq = [char]0x201D

function Refresh-Paragraph([int]$index) {
    # Make a clean duplicate of paragraph $index right after itself (this
    # carries over the run formatting faithfully but drops any stale
    # proof-reading marks the original paragraph was carrying), then
    # delete the stale original so the clean copy takes its place at the
    # same paragraph index.
    $p = $d.Paragraphs($index)
    $p.Range.FormattedText.Copy()
    $afterRange = $d.Range($p.Range.End, $p.Range.End)
    $afterRange.Paste()
    $stale = $d.Paragraphs($index)
    $stale.Range.Delete()
}

# --- "Ls" bullet (paragraph 8): merge " :" and the trailing space that
# follows it into a single " : " run. ---
Refresh-Paragraph 8
$d.Paragraphs(8).Range.Find.Execute(" :" + [char]0x0020, $false, $false, $false, $false, $false, $true, 1, $false, " : ", 2)

# --- "Git commit -m  "..."" bullet (paragraph 9): merge the three runs
# making up the bolded command text into a single run. ---
Refresh-Paragraph 9
$gitCommitText = "Git commit -m  " + $ldq + "descriptive message about the commit" + $rdq
$d.Paragraphs(9).Range.Find.Execute($gitCommitText, $false, $false, $false, $false, $false, $true, 1, $false, $gitCommitText, 2)

# --- "Git push origin master" -> "Git push origin main" ---
$d.Content.Find.Execute("Git push origin master", $false, $false, $false, $false, $false, $true, 1, $false, "Git push origin main", 2)
